$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark rows 35-39 "Completed" (column E) as "Yes", matching the style
#     already used by the existing "Yes" cells above (e.g. E34). ---
$ws.Range("E34").Copy()
$ws.Range("E35:E39").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("E35").Value = "Yes"
$ws.Range("E36").Value = "Yes"
$ws.Range("E37").Value = "Yes"
$ws.Range("E38").Value = "Yes"
$ws.Range("E39").Value = "Yes"

# --- Copy the date style (s="7") used by column D in the existing rows
#     down onto the new rows 40-44 before writing the dates. ---
$ws.Range("D39").Copy()
$ws.Range("D40:D44").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Add the five new "Scheduled Tasks" rows (40-44). ---
$ws.Range("A40").Value = 39
$ws.Range("B40").Value = "Programming Java implementation"
$ws.Range("C40").Value = "Yeqing Liu"
$ws.Range("D40").Value = 43370

$ws.Range("A41").Value = 40
$ws.Range("B41").Value = "Work on Documentation"
$ws.Range("C41").Value = "Timothy Finn"
$ws.Range("D41").Value = 43370

$ws.Range("A42").Value = 41
$ws.Range("B42").Value = "Continue work on Poster for presentation"
$ws.Range("C42").Value = "George Proios"
$ws.Range("D42").Value = 43370

$ws.Range("A43").Value = 42
$ws.Range("B43").Value = "workd on LED and heat map implementation"
$ws.Range("C43").Value = "Josh Francis"
$ws.Range("D43").Value = 43370

$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "Work on integrating LED drivers with Raspberry PI"
$ws.Range("C44").Value = "Michael Douglas"
$ws.Range("D44").Value = 43370

# --- Grow the "Task List" table so the new rows are included. ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E44"))

# --- Update the view: scroll down and select the newly added rows. ---
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1
$ws.Range("A40:E44").Select()
